$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New time-log entry for row 6: set number format to match the other
# Start/End time cells (h:mm), then fill in the times.
$ws.Range("B6").NumberFormat = "h:mm"
$ws.Range("C6").NumberFormat = "h:mm"
$ws.Range("B6").Value = 0.44444444444444442
$ws.Range("C6").Value = 0.5

# Notes column entries. Set in this order so new shared-strings are
# appended matching: "Research factor..." (idx6), "Working up website..."
# (idx7), "Reading more on cluster..." (idx8).
$ws.Range("E6").Value = "Research factor analysis and principal components analysis methods. Exploratory coding work to look at viability of these methods."
$ws.Range("E5").Value = "Working up website to show updates to clustering procedures. "
$ws.Range("E4").Value = "Reading more on cluster methods and beginning work on website. "

# Update the active selection to reflect where the author left off editing.
$ws.Range("E4").Select() | Out-Null
